$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.30599947411284
$ws.Range("C2").Value = 14.91255542715753
$ws.Range("D2").Value = 6.003839878556672
$ws.Range("E2").Value = 11.51856082462967
$ws.Range("F2").Value = 47.21339889380204
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 30.9263124588765
$ws.Range("J2").Value = 9.994734767596624
$ws.Range("M2").Value = 19.22609602362512

$ws.Range("B3").Value = 16.95595404734565
$ws.Range("C3").Value = 14.54505901709155
$ws.Range("D3").Value = 6.007956793367824
$ws.Range("E3").Value = 11.54033031234053
$ws.Range("F3").Value = 46.95490457852198
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 30.84837308728994
$ws.Range("J3").Value = 10.01472407255172
$ws.Range("M3").Value = 19.14821252242033

$ws.Range("B4").Value = 16.74316502880965
$ws.Range("C4").Value = 14.31988694452883
$ws.Range("D4").Value = 6.010771525705141
$ws.Range("E4").Value = 11.55504241591248
$ws.Range("F4").Value = 46.8073777819893
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 30.80705852837171
$ws.Range("J4").Value = 10.02803715559496
$ws.Range("M4").Value = 19.10509357052719

$ws.Range("B5").Value = 16.65712767487054
$ws.Range("C5").Value = 14.22840015692864
$ws.Range("D5").Value = 6.011990849538551
$ws.Range("E5").Value = 11.56137640673815
$ws.Range("F5").Value = 46.7501059617677
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 30.79186855151279
$ws.Range("J5").Value = 10.0337238605642
$ws.Range("M5").Value = 19.08871772725796

$ws.Range("B6").Value = 16.64288610416667
$ws.Range("C6").Value = 14.21322989339798
$ws.Range("D6").Value = 6.012197688344459
$ws.Range("E6").Value = 11.56244862693671
$ws.Range("F6").Value = 46.74076880942214
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 30.78944569143558
$ws.Range("J6").Value = 10.03468393167081
$ws.Range("M6").Value = 19.08607108948771

$ws.Range("B7").Value = 16.7420017796493
$ws.Range("C7").Value = 14.31865181120771
$ws.Range("D7").Value = 6.010787676987774
$ws.Range("E7").Value = 11.55512646654187
$ws.Range("F7").Value = 46.80659382911445
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 30.80684700591389
$ws.Range("J7").Value = 10.02811278948539
$ws.Range("M7").Value = 19.10486786354656

$ws.Range("B8").Value = 17.1849417846218
$ws.Range("C8").Value = 14.78583557908038
$ws.Range("D8").Value = 6.005199924541668
$ws.Range("E8").Value = 11.52578788036137
$ws.Range("F8").Value = 47.12196869601308
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 30.89808278264798
$ws.Range("J8").Value = 10.00141138196097
$ws.Range("M8").Value = 19.19827428401732

$ws.Range("B9").Value = 18.06457266199416
$ws.Range("C9").Value = 15.69922619141932
$ws.Range("D9").Value = 5.996513010070069
$ws.Range("E9").Value = 11.4789165512848
$ws.Range("F9").Value = 47.82755166858816
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 31.12876371614958
$ws.Range("J9").Value = 9.957293919893376
$ws.Range("M9").Value = 19.4181035575444

$ws.Range("B10").Value = 18.70981898364641
$ws.Range("C10").Value = 16.36037258554948
$ws.Range("D10").Value = 5.991507554871694
$ws.Range("E10").Value = 11.45096006215017
$ws.Range("F10").Value = 48.39656286102626
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 31.3295344892484
$ws.Range("J10").Value = 9.929900547935816
$ws.Range("M10").Value = 19.6009985650148

$ws.Range("B11").Value = 19.00159310800284
$ws.Range("C11").Value = 16.65740348452536
$ws.Range("D11").Value = 5.989527976062719
$ws.Range("E11").Value = 11.43964472229599
$ws.Range("F11").Value = 48.66579617806083
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 31.42757284444646
$ws.Range("J11").Value = 9.918527386381012
$ws.Range("M11").Value = 19.6886131602659

$ws.Range("B12").Value = 19.11170809485709
$ws.Range("C12").Value = 16.76922498163942
$ws.Range("D12").Value = 5.988821019681295
$ws.Range("E12").Value = 11.4355611752358
$ws.Range("F12").Value = 48.76918395575424
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 31.46565123484247
$ws.Range("J12").Value = 9.914377070873229
$ws.Range("M12").Value = 19.72240384149387

$ws.Range("B13").Value = 19.08801133884682
$ws.Range("C13").Value = 16.74517330884036
$ws.Range("D13").Value = 5.988971379325928
$ws.Range("E13").Value = 11.43643169111445
$ws.Range("F13").Value = 48.74685469515209
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 31.45740816207791
$ws.Range("J13").Value = 9.915263957717638
$ws.Range("M13").Value = 19.71509951321244

$ws.Range("B14").Value = 19.01066049279146
$ws.Range("C14").Value = 16.66661695046903
$ws.Range("D14").Value = 5.989468959786434
$ws.Range("E14").Value = 11.43930473333618
$ws.Range("F14").Value = 48.6742735412582
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 31.43068653791804
$ws.Range("J14").Value = 9.918182802027642
$ws.Range("M14").Value = 19.69138098503816

$ws.Range("B15").Value = 18.96322872728141
$ws.Range("C15").Value = 16.61840981510671
$ws.Range("D15").Value = 5.989779295772767
$ws.Range("E15").Value = 11.44109076400641
$ws.Range("F15").Value = 48.63000056772413
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 31.41444257713562
$ws.Range("J15").Value = 9.919991052480345
$ws.Range("M15").Value = 19.67693185572198

$ws.Range("B16").Value = 18.69070505486521
$ws.Range("C16").Value = 16.34087529875236
$ws.Range("D16").Value = 5.991642901097876
$ws.Range("E16").Value = 11.45172773597124
$ws.Range("F16").Value = 48.37917190235643
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 31.32326150249447
$ws.Range("J16").Value = 9.930665705895873
$ws.Range("M16").Value = 19.59535973394322

$ws.Range("B17").Value = 18.52298325354576
$ws.Range("C17").Value = 16.16957262883157
$ws.Range("D17").Value = 5.992862266029611
$ws.Range("E17").Value = 11.45861209519088
$ws.Range("F17").Value = 48.22791631494118
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 31.26903528289908
$ws.Range("J17").Value = 9.937492956927468
$ws.Range("M17").Value = 19.54643336957711

$ws.Range("B18").Value = 18.4263560940541
$ws.Range("C18").Value = 16.07069966837128
$ws.Range("D18").Value = 5.993591613551974
$ws.Range("E18").Value = 11.46270379792838
$ws.Range("F18").Value = 48.1418995293953
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 31.23847782757118
$ws.Range("J18").Value = 9.941522241404535
$ws.Range("M18").Value = 19.51870913882529

$ws.Range("B19").Value = 18.39361647727131
$ws.Range("C19").Value = 16.03716772240432
$ws.Range("D19").Value = 5.993843370533387
$ws.Range("E19").Value = 11.46411185881605
$ws.Range("F19").Value = 48.11294602146576
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 31.22824045998492
$ws.Range("J19").Value = 9.942904081347201
$ws.Range("M19").Value = 19.50939444626687

$ws.Range("B20").Value = 18.54085478268408
$ws.Range("C20").Value = 16.18784462344628
$ws.Range("D20").Value = 5.992729565310338
$ws.Range("E20").Value = 11.45786558397353
$ws.Range("F20").Value = 48.24391657171353
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 31.27474241846818
$ws.Range("J20").Value = 9.936755583489532
$ws.Range("M20").Value = 19.55159867014394

$ws.Range("B21").Value = 19.03339138430209
$ws.Range("C21").Value = 16.68970962904896
$ws.Range("D21").Value = 5.98932165115705
$ws.Range("E21").Value = 11.43845538975019
$ws.Range("F21").Value = 48.69555392428742
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 31.43850954757136
$ws.Range("J21").Value = 9.917321220822316
$ws.Range("M21").Value = 19.69833123662745

$ws.Range("B22").Value = 19.3530637757601
$ws.Range("C22").Value = 17.01382421664796
$ws.Range("D22").Value = 5.987343035931986
$ws.Range("E22").Value = 11.42694304197228
$ws.Range("F22").Value = 48.9990597743534
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 31.55109171891118
$ws.Range("J22").Value = 9.905531661524863
$ws.Range("M22").Value = 19.79779122209792

$ws.Range("B23").Value = 19.18269091916511
$ws.Range("C23").Value = 16.84123125621254
$ws.Range("D23").Value = 5.988376340561331
$ws.Range("E23").Value = 11.43298014247885
$ws.Range("F23").Value = 48.83633053885017
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 31.49050063726097
$ws.Range("J23").Value = 9.91174054018974
$ws.Range("M23").Value = 19.74438915839003

$ws.Range("B24").Value = 18.53277568318595
$ws.Range("C24").Value = 16.17958505629757
$ws.Range("D24").Value = 5.99278947109239
$ws.Range("E24").Value = 11.45820266493663
$ws.Range("F24").Value = 48.23667991688506
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 31.27216029697536
$ws.Range("J24").Value = 9.937088625437298
$ws.Range("M24").Value = 19.54926217715532

$ws.Range("B25").Value = 17.82624497785569
$ws.Range("C25").Value = 15.45330496409028
$ws.Range("D25").Value = 5.99862075213244
$ws.Range("E25").Value = 11.49045710355245
$ws.Range("F25").Value = 47.62758600425985
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 31.06083744745467
$ws.Range("J25").Value = 9.968346919989317
$ws.Range("M25").Value = 19.35480307913805

